# Fix misspelling: 라이타니엔 -> 라이타니아 (Laitania) in the Korean (ko_KR) story column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9   = '[name="살카즈 전사"]  라이타니아 각 도시에 우뚝 솟은 탑에서 평민들 깔보고 사는 대귀족 놈 중 한 놈만 털어도, 우리 모두가 배불리 먹고 마실 수 있을 텐데.
'
    11  = '[name="살카즈 전사"]  하지만 그 오만함 덕분에 라이타니아의 추악한 침략 본능이 억제될 수 있었지.
'
    45  = '[name="그레이스롯"]  ……게다가 나도 진상이 어떠한들, ''리유니온이 다시 나타나 라이타니아의 마을을 괴멸시키는'' 그런 일은……
'
    63  = '[name="붙잡힌 감염자"]  웃기지 마…… 난 라이타니아인이야, 고등 교육을 받은 라이타니아인이란 말이다. 그런데 어떻게……
'
    95  = '[name="할아버지"]  우린 라이타니아을, 라이타니아의 법을, 이동도시의 엔진에서 나는 굉음을, 결국 꺾지 못한 게다……
'
    106 = '[name="스즈란"]  비더만은 라이타니아의 한 재앙정보전달자의 이름이잖아요. 그 사람은 분명 화재에 희생되었던 게……
'
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
